# Resize the single inline text-box shape ("文本框 2") that holds the
# "git cherry-pick <commitHash> //commit标签" example, shrinking its
# width from 3372416 EMU (265.55pt) to 3182293 EMU (250.57pt).
# Height (253497 EMU / 19.95pt) is unchanged.

$d = $word.ActiveDocument

$shp = $d.InlineShapes.Item(1)

# Target width in points, computed from the target EMU extent
# (1 inch = 914400 EMU = 72 points):
$shp.Width  = [double](3182293 / 914400.0 * 72)
$shp.Height = [double](253497 / 914400.0 * 72)
